# Applies the "Updated symbol list" commit: refreshed prices/volumes for
# several coins, a 3-way row rotation among KickToken/BKEXToken/CEJI
# (rows 41-43), and a couple of label tweaks (E47).
#
# Price cells in column D are stored as *text* (not numbers) in the
# original workbook, so NumberFormat is forced to "@" (Text) before each
# assignment - otherwise Excel would auto-coerce the numeric-looking
# string into a Number cell, which would not match the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.91"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.395"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06005"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8119"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9544"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07407"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03055"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09423"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.003"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001589"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04793"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005872"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006243"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005063"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009887"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.701"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.402"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1340"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002461"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04016"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1074"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002721"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003020"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005848"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005265"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8013"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.02255"
